# Apply updated cryptocurrency price/volume figures per the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.079.65"
$ws.Range("E2").Value = "  -0.24%  "

$ws.Range("D3").Value = "1.622.89"
$ws.Range("E3").Value = "  -1.01%  "

$ws.Range("E4").Value = "  -0.06%  "

$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "214.24"
$c.ClearFormats()
$ws.Range("E5").Value = "  -1.13%  "

$ws.Range("E6").Value = "  -1.13%  "

$ws.Range("E7").Value = "  -0.05%  "

$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = "0.0629"
$c.ClearFormats()
$ws.Range("E8").Value = "  +0.59%  "

$ws.Range("E9").Value = "  -1.64%  "

$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "20.00"
$c.ClearFormats()
$ws.Range("E10").Value = "  +0.19%  "

$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "0.0846"
$c.ClearFormats()
$ws.Range("E11").Value = "  -0.07%  "

$ws.Range("D12").Value = "1.852.23"
$ws.Range("E12").Value = "  -0.89%  "

$ws.Range("D13").Value = "1.622.98"
$ws.Range("E13").Value = "  -0.87%  "

$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "4.13"
$c.ClearFormats()
$ws.Range("E14").Value = "  +0.00%  "

$ws.Range("E15").Value = "  -0.43%  "

$ws.Range("B16").Value = "WrappedBTC"
$ws.Range("C16").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D16").Value = "27.057.50"
$ws.Range("E16").Value = "  -0.33%  "

$ws.Range("B17").Value = "Litecoin"
$ws.Range("C17").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = "64.46"
$c.ClearFormats()
$ws.Range("E17").Value = "  -3.55%  "

$ws.Range("D18").Value = "0.0₃0739"
$ws.Range("E18").Value = "  -0.10%  "

$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "212.98"
$c.ClearFormats()
$ws.Range("E19").Value = "  -1.97%  "

$ws.Range("E20").Value = "  -0.09%  "

$ws.Range("E21").Value = "  -2.05%  "

$ws.Range("E22").Value = "  -1.50%  "

$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "2.35"
$c.ClearFormats()
$ws.Range("E23").Value = "  -7.12%  "

$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "9.05"
$c.ClearFormats()
$ws.Range("E24").Value = "  -0.88%  "

$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "147.86"
$c.ClearFormats()
$ws.Range("E25").Value = "  +0.47%  "

$ws.Range("E26").Value = "  +0.05%  "

$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "7.36"
$c.ClearFormats()
$ws.Range("E27").Value = "  -0.88%  "

$ws.Range("E28").Value = "  -3.10%  "

$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "15.53"
$c.ClearFormats()
$ws.Range("E29").Value = "  -0.97%  "

$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "0.0510"
$c.ClearFormats()
$ws.Range("E30").Value = "  +0.22%  "

$ws.Range("E31").Value = "  -0.90%  "

$ws.Range("E32").Value = "  -1.19%  "

$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "0.726"
$c.ClearFormats()
$ws.Range("E33").Value = "  +33.53%  "

$ws.Range("E34").Value = "  -0.76%  "

$ws.Range("D35").Value = "1.360.39"
$ws.Range("E35").Value = "  +4.10%  "

$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "1.56"
$c.ClearFormats()
$ws.Range("E36").Value = "  -0.15%  "

$ws.Range("E37").Value = "  -0.54%  "

$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "0.0176"
$c.ClearFormats()
$ws.Range("E38").Value = "  +0.41%  "

$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "0.842"
$c.ClearFormats()
$ws.Range("E39").Value = "  -1.84%  "

$ws.Range("E40").Value = "  -0.10%  "

$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "0.802"
$c.ClearFormats()
$ws.Range("E41").Value = "  -1.00%  "

$ws.Range("E42").Value = "  +0.61%  "

$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "64.41"
$c.ClearFormats()
$ws.Range("E43").Value = "  +4.28%  "

$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "5.35"
$c.ClearFormats()
$ws.Range("E44").Value = "  +0.90%  "

$ws.Range("D45").Value = "1.762.42"
$ws.Range("E45").Value = "  -0.92%  "

$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "1.65"
$c.ClearFormats()
$ws.Range("E46").Value = "  +3.35%  "

$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "89.86"
$c.ClearFormats()
$ws.Range("E47").Value = "  -1.85%  "

$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "0.867"
$c.ClearFormats()
$ws.Range("E48").Value = "  +30.22%  "

$ws.Range("E49").Value = "  -1.72%  "

$ws.Range("E50").Value = "  +4.83%  "

$ws.Range("E51").Value = "  +0.17%  "
